$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-script ("Description", column G) and validate-block
# ("Expected Behaviour", column H) cell contents for the affected test cases.
# These correspond to shared-string entries that were edited/reordered in the
# original OOXML; here we simply assign the final text directly to each cell.
$ws.Range("G3").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0749_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`nwifi_Mode(OFF);`nwait(2);`npress_Key(Home);`nlaunch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);`nvalidate5;`nwifi_Mode(ON);`nwait(20);`npress_Key(Home);"
$ws.Range("G4").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0750_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`nwifi_Mode(OFF);`nwait(2);`npress_Key(Home);`nlaunch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);`nvalidate5;`nwifi_Mode(ON);`nwait(20);`npress_Key(Home);"
$ws.Range("H4").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Network`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0750`n};`nvalidate4`n{`nvalidate_Result=Connected`n};`nvalidate5`n{`nvalidate_doesNotContain=Disconnected`n};"
$ws.Range("G5").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0751_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`nwifi_Mode(OFF);`nwait(2);`npress_Key(Home);`nlaunch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);`nvalidate5;`nwifi_Mode(ON);`nwait(20);`npress_Key(Home);"
$ws.Range("G6").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0752_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`nwifi_Mode(OFF);`nwait(2);`npress_Key(Home);`nlaunch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);`nvalidate5;`nwifi_Mode(ON);`nwait(20);`npress_Key(Home);"
$ws.Range("G7").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0755_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`nwifi_Mode(OFF);`nwait(60);`npress_Key(Home);`nlaunch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);`nvalidate5;`nwifi_Mode(ON);`nwait(20);`npress_Key(Home);"

# Align D8 formatting with the other cells in column D (D2:D7), which were
# using a slightly different cell style (removes the stray top/bottom border
# variant and extra fill flag that only D8 had).
$ws.Range("D2").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection shown when the workbook is opened.
$ws.Range("A2").Select()
